$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before D (shifts existing D:K data to F:M)
$ws.Range("D:E").EntireColumn.Insert()

# Copy number formatting/style from the (now-shifted) data columns F:G back
# onto the new D:E columns for each contiguous data block, so the new
# columns pick up the correct date/number styles instead of the default
# left-column style.
$ws.Range("F7:G35").Copy()
$ws.Range("D7").PasteSpecial(-4122)

$ws.Range("F38:G77").Copy()
$ws.Range("D38").PasteSpecial(-4122)

$ws.Range("F80:G102").Copy()
$ws.Range("D80").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Populate the two new quarters of data in columns D (most recent) and E
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 274100
$ws.Range("E8").Value = 274600
$ws.Range("D9").Value = 215600
$ws.Range("E9").Value = 215700
$ws.Range("D10").Value = 58500
$ws.Range("E10").Value = 58900
$ws.Range("D12").Value = 'NA'
$ws.Range("E12").Value = 'NA'
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 32800
$ws.Range("E15").Value = 30600
$ws.Range("D17").Value = 284100
$ws.Range("E17").Value = 278300
$ws.Range("D18").Value = -10000
$ws.Range("E18").Value = -3700
$ws.Range("D20").Value = 1300
$ws.Range("E20").Value = 800
$ws.Range("D21").Value = 24200
$ws.Range("E21").Value = 27600
$ws.Range("D22").Value = 5000
$ws.Range("E22").Value = 4900
$ws.Range("D23").Value = -13600
$ws.Range("E23").Value = -7900
$ws.Range("D24").Value = 700
$ws.Range("E24").Value = 2000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -14300
$ws.Range("E26").Value = -9800
$ws.Range("D27").Value = -14300
$ws.Range("E27").Value = -9800
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 5800
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -1300
$ws.Range("E32").Value = -800
$ws.Range("D33").Value = -14300
$ws.Range("E33").Value = -4000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -14300
$ws.Range("E35").Value = -4000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 19300
$ws.Range("E41").Value = 36300
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 293100
$ws.Range("E43").Value = 306400
$ws.Range("D44").Value = 209400
$ws.Range("E44").Value = 210800
$ws.Range("D45").Value = 12200
$ws.Range("E45").Value = 12100
$ws.Range("D46").Value = 534000
$ws.Range("E46").Value = 565600
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 540400
$ws.Range("E48").Value = 544700
$ws.Range("D49").Value = 902300
$ws.Range("E49").Value = 904600
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 27000
$ws.Range("E52").Value = 29900
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 2003800
$ws.Range("E54").Value = 2044800
$ws.Range("D57").Value = 77500
$ws.Range("E57").Value = 78600
$ws.Range("D58").Value = 25600
$ws.Range("E58").Value = 25500
$ws.Range("D59").Value = 78000
$ws.Range("E59").Value = 81700
$ws.Range("D60").Value = 181000
$ws.Range("E60").Value = 185900
$ws.Range("D61").Value = 306200
$ws.Range("E61").Value = 328900
$ws.Range("D62").Value = 76800
$ws.Range("E62").Value = 80400
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 564100
$ws.Range("E66").Value = 595100
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 1029500
$ws.Range("E72").Value = 1043900
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1439800
$ws.Range("E76").Value = 1449700
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -14300
$ws.Range("E81").Value = -4000
$ws.Range("D83").Value = 32800
$ws.Range("E83").Value = 30600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 23100
$ws.Range("E89").Value = 33200
$ws.Range("D91").Value = -16700
$ws.Range("E91").Value = -33000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -14600
$ws.Range("E94").Value = -29100
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -24600
$ws.Range("E100").Value = 2400
$ws.Range("D101").Value = -800
$ws.Range("E101").Value = 700
$ws.Range("D102").Value = -16900
$ws.Range("E102").Value = 7200
